$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 38.04655766666666
$ws.Range("H2").Value = 114.139673
$ws.Range("I2").Value = 0.8090698722086991
$ws.Range("J2").Value = 0.8090698722086992
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.722664999999999
$ws.Range("N2").Value = 17.167995
$ws.Range("O2").Value = 0.09021166427595352
$ws.Range("P2").Value = 0.09021166427595351
$ws.Range("Q2").Value = 217.727703929515
$ws.Range("R2").Value = 1959.549335365635
$ws.Range("S2").Value = 0.07298753968747979
$ws.Range("T2").Value = 0.07298753968747977

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 38.04655766666666
$ws.Range("H3").Value = 114.139673
$ws.Range("I3").Value = 0.8090698722086991
$ws.Range("J3").Value = 0.8090698722086992
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 41.286995
$ws.Range("N3").Value = 123.860985
$ws.Range("O3").Value = 0.6508451100847196
$ws.Range("P3").Value = 0.6508451100847196
$ws.Range("Q3").Value = 1570.828036150878
$ws.Range("R3").Value = 14137.4523253579
$ws.Range("S3").Value = 0.5265791700439008
$ws.Range("T3").Value = 0.5265791700439009

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 38.04655766666666
$ws.Range("H4").Value = 114.139673
$ws.Range("I4").Value = 0.8090698722086991
$ws.Range("J4").Value = 0.8090698722086992
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06212466666666667
$ws.Range("N4").Value = 0.186374
$ws.Range("O4").Value = 0.0009793286122093212
$ws.Range("P4").Value = 0.000979328612209321
$ws.Range("Q4").Value = 2.363629712855778
$ws.Range("R4").Value = 21.272667415702
$ws.Range("S4").Value = 0.0007923452751305182
$ws.Range("T4").Value = 0.000792345275130518

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 38.04655766666666
$ws.Range("H5").Value = 114.139673
$ws.Range("I5").Value = 0.8090698722086991
$ws.Range("J5").Value = 0.8090698722086992
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.36419166666667
$ws.Range("N5").Value = 49.092575
$ws.Range("O5").Value = 0.2579638970271176
$ws.Range("P5").Value = 0.2579638970271176
$ws.Range("Q5").Value = 622.6011619142195
$ws.Range("R5").Value = 5603.410457227976
$ws.Range("S5").Value = 0.208710817202188
$ws.Range("T5").Value = 0.2087108172021881

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.03439113957782537
$ws.Range("J6").Value = 0.03439113957782537
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.722664999999999
$ws.Range("N6").Value = 17.167995
$ws.Range("O6").Value = 0.09021166427595352
$ws.Range("P6").Value = 0.09021166427595351
$ws.Range("Q6").Value = 9.254953265480001
$ws.Range("R6").Value = 83.29457938932
$ws.Range("S6").Value = 0.00310248193766224
$ws.Range("T6").Value = 0.00310248193766224

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.03439113957782537
$ws.Range("J7").Value = 0.03439113957782537
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 41.286995
$ws.Range("N7").Value = 123.860985
$ws.Range("O7").Value = 0.6508451100847196
$ws.Range("P7").Value = 0.6508451100847196
$ws.Range("Q7").Value = 66.77119999110667
$ws.Range("R7").Value = 600.9407999199601
$ws.Range("S7").Value = 0.02238330502446871
$ws.Range("T7").Value = 0.02238330502446871

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.03439113957782537
$ws.Range("J8").Value = 0.03439113957782537
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06212466666666667
$ws.Range("N8").Value = 0.186374
$ws.Range("O8").Value = 0.0009793286122093212
$ws.Range("P8").Value = 0.000979328612209321
$ws.Range("Q8").Value = 0.1004708272515556
$ws.Range("R8").Value = 0.9042374452640002
$ws.Range("S8").Value = [double]"3.368022699504878e-05"
$ws.Range("T8").Value = [double]"3.368022699504878e-05"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.03439113957782537
$ws.Range("J9").Value = 0.03439113957782537
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.36419166666667
$ws.Range("N9").Value = 49.092575
$ws.Range("O9").Value = 0.2579638970271176
$ws.Range("P9").Value = 0.2579638970271176
$ws.Range("Q9").Value = 26.46491260668889
$ws.Range("R9").Value = 238.1842134602001
$ws.Range("S9").Value = 0.008871672388699372
$ws.Range("T9").Value = 0.008871672388699372

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Bmpr1a"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.361255
$ws.Range("H10").Value = 22.083765
$ws.Range("I10").Value = 0.1565389882134754
$ws.Range("J10").Value = 0.1565389882134754
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.722664999999999
$ws.Range("N10").Value = 17.167995
$ws.Range("O10").Value = 0.09021166427595352
$ws.Range("P10").Value = 0.09021166427595351
$ws.Range("Q10").Value = 42.12599634457499
$ws.Range("R10").Value = 379.133967101175
$ws.Range("S10").Value = 0.01412164265081149
$ws.Range("T10").Value = 0.01412164265081149

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Bmp6"
$ws.Range("C11").Value = "Bmpr1a"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.361255
$ws.Range("H11").Value = 22.083765
$ws.Range("I11").Value = 0.1565389882134754
$ws.Range("J11").Value = 0.1565389882134754
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 41.286995
$ws.Range("N11").Value = 123.860985
$ws.Range("O11").Value = 0.6508451100847196
$ws.Range("P11").Value = 0.6508451100847196
$ws.Range("Q11").Value = 303.924098378725
$ws.Range("R11").Value = 2735.316885408525
$ws.Range("S11").Value = 0.10188263501635
$ws.Range("T11").Value = 0.10188263501635

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Bmp6"
$ws.Range("C12").Value = "Bmpr1a"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.361255
$ws.Range("H12").Value = 22.083765
$ws.Range("I12").Value = 0.1565389882134754
$ws.Range("J12").Value = 0.1565389882134754
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.06212466666666667
$ws.Range("N12").Value = 0.186374
$ws.Range("O12").Value = 0.0009793286122093212
$ws.Range("P12").Value = 0.000979328612209321
$ws.Range("Q12").Value = 0.4573155131233334
$ws.Range("R12").Value = 4.11583961811
$ws.Range("S12").Value = 0.0001533031100837542
$ws.Range("T12").Value = 0.0001533031100837542

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Bmp6"
$ws.Range("C13").Value = "Bmpr1a"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.361255
$ws.Range("H13").Value = 22.083765
$ws.Range("I13").Value = 0.1565389882134754
$ws.Range("J13").Value = 0.1565389882134754
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 16.36419166666667
$ws.Range("N13").Value = 49.092575
$ws.Range("O13").Value = 0.2579638970271176
$ws.Range("P13").Value = 0.2579638970271176
$ws.Range("Q13").Value = 120.4609877272083
$ws.Range("R13").Value = 1084.148889544875
$ws.Range("S13").Value = 0.04038140743623016
$ws.Range("T13").Value = 0.04038140743623016
